$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Binning" row (row 6) values from "[]" to "[-1,60,80,100]"
# across all data columns (B:M).
$ws.Range("B6:M6").Value = "[-1,60,80,100]"

# Reflect the new selection left by the editing session: C6:M6, active cell C6.
$ws.Range("C6:M6").Select() | Out-Null
